# Update the "K" column (column G) values on Sheet1 for rows 2-36.
# These values were regenerated from strikeout box-score data
# (commit: "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newK = @{
    2  = 2
    3  = 3
    4  = 2
    5  = 1
    6  = 1
    7  = 1
    8  = 1
    9  = 2
    10 = 3
    11 = 7
    12 = 2
    13 = 3
    14 = 5
    15 = 1
    16 = 2
    17 = 5
    18 = 1
    19 = 7
    20 = 2
    21 = 5
    22 = 2
    23 = 0
    24 = 3
    25 = 1
    26 = 0
    27 = 9
    28 = 3
    29 = 3
    30 = 1
    31 = 3
    32 = 4
    33 = 6
    34 = 2
    35 = 1
    36 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
